# Generate Report for Handoff
# Appends two new file records (9aa33e59-7ea7-44d1-be7a-8f7773df1308 and
# c5da8a89-02de-45f7-961e-d5acc6e56cc9) to the Overview, zh-cn and de-de
# sheets, mirroring the existing rows for the other files already present.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTest/oltest/blob/19989c9daf74c7635e5a196b116948416e9fc42e/e2e"
$zhHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3fac1f8bc2e9b1c3265c31c2b8df49cbe3818ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht"
$deHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf8d7e1ed582dd3db001823dcf594c161ad3eabd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht"

$file1 = "9aa33e59-7ea7-44d1-be7a-8f7773df1308"
$file1Hash = "a402ae9715023b7e114360eecef391ad1ab09a9f"
$file1ZhDate = "2016-03-09 06:21:29"
$file1DeDate = "2016-03-09 06:21:39"
$file1OverviewDate = "2016-21-09 06:21:39"

$file2 = "c5da8a89-02de-45f7-961e-d5acc6e56cc9"
$file2Hash = "6d9d1218c2c5696cad9d730aab67f406db4adfbb"
$file2ZhDate = $file1ZhDate
$file2DeDate = $file1DeDate
$file2OverviewDate = $file1OverviewDate

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = $file1OverviewDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$ghBase/$file1.md", "", "", "$file1.md")

$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = $file2OverviewDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "$ghBase/$file2.md", "", "", "$file2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Handoff Reason |
#   Dependency From | Error Detail
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("E4").Value = $file1ZhDate
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$ghBase/$file1.md", "", "", "$file1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), "$ghBase/$file1.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "$zhHtBase/$file1.$file1Hash.zh-cn.xlf", "", "", "$file1.$file1Hash.zh-cn.xlf")

$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("E5").Value = $file2ZhDate
$wsZh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "$ghBase/$file2.md", "", "", "$file2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), "$ghBase/$file2.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "$zhHtBase/$file2.$file2Hash.zh-cn.xlf", "", "", "$file2.$file2Hash.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("E4").Value = $file1DeDate
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$ghBase/$file1.md", "", "", "$file1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), "$ghBase/$file1.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "$deHtBase/$file1.$file1Hash.de-de.xlf", "", "", "$file1.$file1Hash.de-de.xlf")

$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("E5").Value = $file2DeDate
$wsDe.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "$ghBase/$file2.md", "", "", "$file2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), "$ghBase/$file2.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "$deHtBase/$file2.$file2Hash.de-de.xlf", "", "", "$file2.$file2Hash.de-de.xlf")

Write-Host "Added handoff rows for $file1 and $file2 to Overview, zh-cn and de-de sheets."
